# Apply the "time period (unit:s)" column edit to the "new LFC function list" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("new LFC function list")
$ws.Activate()

# New header for column N
$ws.Range("N1").Value = "time period(unit:s)"

# Fill in the time-period values (column N) for the data rows.
# Rows using 120 seconds
$ws.Range("N4:N16").Value = 120
$ws.Range("N31:N61").Value = 120
$ws.Range("N66:N72").Value = 120

# Rows using 20 seconds
$ws.Range("N19:N27").Value = 20
$ws.Range("N30").Value = 20
$ws.Range("N63:N64").Value = 20

# Widen/narrow a few columns to match the final layout
$ws.Columns.Item(4).ColumnWidth = 10.833333333333334
$ws.Columns.Item(7).ColumnWidth = 34.5
$ws.Columns.Item(13).ColumnWidth = 25.666666666666668
$ws.Columns.Item(14).ColumnWidth = 14.5

# Restore the view scroll position / active selection as left by the author
$ws.Application.ActiveWindow.ScrollRow = 59
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("O75").Select()
